# Adding mission modifiers for Legendary dragons
# - Insert 5 new rows into the missionDragonModifiersDefinitions table
#   (Table13303132) on the "missions" sheet, for dragons: dragon_electric,
#   dragon_helicopter, dragon_hedgehog, dragon_ice, dragon_dino.
# - Everything below shifts down by 5 rows (other tables, plain data block).
# - Update dependent table refs (autoFilter/ref) to match new row numbers.
# - Update sheet selection / active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("missions")

# --- 1. Make room: insert 5 blank rows right after the last existing
#        dragon-modifier row (94), shifting every row below down by 5.
$ws.Rows.Item(95).Resize(5).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown) | Out-Null

# --- 2. Resize the dependent ListObjects so their XML ref/autoFilter match
#        the new layout (the dragon-modifiers table grows by 5 rows; the
#        difficulty- and other-modifiers tables merely shift down by 5).
#        This MUST happen before any of the new rows are filled in below:
#        ListObject.Range is cached at the pre-insert coordinates, and
#        Resize() re-reads the (now shifted-into-place) header text from
#        there, so writing new data first would clobber the headers it
#        reads.
$loDragon = $ws.ListObjects.Item("Table13303132")
$loDragon.Resize($ws.Range("B81:E99")) | Out-Null

$loDifficulty = $ws.ListObjects.Item("Table1330313234")
$loDifficulty.Resize($ws.Range("B103:E106")) | Out-Null

$loOther = $ws.ListObjects.Item("Table133031323435")
$loOther.Resize($ws.Range("B110:D111")) | Out-Null

# --- 3. Copy the formatting of the last existing data row (94) across the
#        5 new rows (95:99) so styles/number formats match the table.
$ws.Range("B94:E94").Copy() | Out-Null
$ws.Range("B95:E99").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# --- 4. Fill in the new dragon modifier rows.
$skus = @("dragon_electric", "dragon_helicopter", "dragon_hedgehog", "dragon_ice", "dragon_dino")
for ($i = 0; $i -lt $skus.Length; $i++) {
    $r = 95 + $i
    $ws.Cells.Item($r, 2).Value2 = "<Definition>"
    $ws.Cells.Item($r, 3).Value2 = $skus[$i]
    $ws.Cells.Item($r, 4).Value2 = 8
    $ws.Cells.Item($r, 5).Value2 = 90
}

# --- 5. Update view state: "missions" becomes the active/selected sheet
#        (was "tournaments"), with a fresh selection/scroll position.
$ws.Activate()
$ws.Range("G95").Select() | Out-Null
